$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.780.77'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '2.908.04'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.05%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '568.69'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.03%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '144.64'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("D9").Value = '2.906.61'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E13").Value = '  +1.08%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '32.91'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.29%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.125'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '3.389.51'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '62.704.35'
$ws.Range("E17").Value = '  +1.47%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.64'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '2.904.59'
$ws.Range("E19").Value = '  -0.05%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '429.37'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.42%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.11'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").Value = '  -0.09%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.89'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.69%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '78.87'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.74%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '11.96'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '10.11'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.98%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("E29").Value = '  +5.46%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -2.35%  '
$ws.Range("E32").Value = '  -3.44%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '25.88'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -2.27%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.955'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("E37").Value = '  -1.81%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.39%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '48.71'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("E41").Value = '  -1.52%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '41.12'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +5.37%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '8.09'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("D45").Value = '2.715.87'
$ws.Range("E45").Value = '  +0.62%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0340'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.19%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '133.23'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.40%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '359.25'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +5.17%  '
$ws.Range("E50").Value = '  +14.48%  '
$ws.Range("E51").Value = '  -0.38%  '
